$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "123/2022 Otsikko"
$ws.Range("A3").Select()
